$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 72: add the missing PM ("nr") marker and the trig-point note ---
$ws.Range("J72").Value = "nr"
$ws.Range("R72").Value = "Star Disk"

# --- Row 73: correct the date (20 Nov 2023 -> 20 Oct 2023) ---
$ws.Range("A73").Value = 45219

# --- Insert two new rows after row 73 (pushes old blank row74 -> 76, ---
# --- old totals row75 -> 77, old trailer row78 -> 80) ---
$ws.Rows.Item(74).Insert()
$ws.Rows.Item(74).Insert()

# --- Fill in the new trip on row 74 ---
$ws.Range("A74").Value = 45252
$ws.Range("B74").Value = "Sawley Junction"
$ws.Range("C74").Value = "Long Eaton"
$ws.Range("D74").Value = $ws.Range("D73").Value2
$ws.Range("E74").Value = 5.37
$ws.Range("F74").Formula = "=E74*0.6213712"
$ws.Range("G74").Value = 0.02613425925925926
$ws.Range("H74").Formula = "=G74/F74"
$ws.Range("I74").Value = 1
$ws.Range("J74").Value = "nr"
$ws.Range("K74").Value = 1
$ws.Range("O74").Value = 1
$ws.Range("P74").Value = "Up to trent lock and back. Sloe Gin Cider"
$ws.Range("Q74").Formula = "=SUM(I74:O74)*F74"

# --- Extend the totals row (now row 77) to include the new row 74 ---
$ws.Range("F77").Formula = "=SUM(F8:F74)"

# --- Column P width ---
$ws.Columns.Item(16).ColumnWidth = 33

# --- Selection / scroll tweaks ---
$null = $excel.ActiveWindow.ScrollColumn
$ws.Application.ActiveWindow.ScrollColumn = 3
$null = $ws.Range("S74").Select()
